$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update column C (Förändrad) for rows 2 through 10 from 45170 to 45174
for ($row = 2; $row -le 10; $row++) {
    $ws.Cells.Item($row, 3).Value = 45174
}
